$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (shifts existing rows 18-23 down to 19-24)
$ws.Rows("18:18").Insert()

# Populate the new row 18 with the new data record
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44449
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100108
$ws.Range("H18").Value = "Tropicales y subtropicales"
$ws.Range("I18").Value = 100108003
$ws.Range("J18").Value = "Maracuyá"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 20
$ws.Range("N18").Value = 38000
$ws.Range("O18").Value = 38000
$ws.Range("P18").Value = 38000
$ws.Range("Q18").Value = "$/caja 18 kilos"
$ws.Range("R18").Value = "Región de Arica y Parinacota"
$ws.Range("S18").Value = 2111
$ws.Range("T18").Value = 18
